# Update Name of Algo
# Apply updated numeric results (column C/D/A values) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = -13.41370000000001
$ws.Range("C4").Value  = -14.49680000000002
$ws.Range("D6").Value  = -8.048500000000001
$ws.Range("C7").Value  = -11.43469999999999
$ws.Range("D7").Value  = -7.483099999999995
$ws.Range("C8").Value  = -12.22599999999998
$ws.Range("D8").Value  = -8.095100000000006
$ws.Range("A11").Value = -21.89910000000001
$ws.Range("A12").Value = -20.83240000000002
$ws.Range("C12").Value = -11.2419
$ws.Range("C14").Value = -12.0213
$ws.Range("A15").Value = -21.17630000000002
$ws.Range("D19").Value = -8.643999999999991
$ws.Range("D21").Value = -7.877099999999996
$ws.Range("C22").Value = -11.17799999999999
$ws.Range("D24").Value = -7.630499999999996
$ws.Range("D25").Value = -7.572299999999999
